# --- Commit: "Sat, Jul 04, 2020  8:05:43 AM" ---
# 1) The comparison table on slide 5 switches from the deck's custom
#    "Table_0" style to the built-in "Medium Style 2 - Accent 1" style.
$p  = $ppt.ActivePresentation
$s5 = $p.Slides.Item(5)
$tblShape = $s5.Shapes.Item(2)
$tblShape.Table.ApplyStyle("{5DC53750-7BD0-4BDC-9EB0-408EFA7932EE}")

# 2) The presentation's theme colours are swapped: the design applied to
#    the slide master changes from the "Integral" / "Red Violet" palette
#    to the stock "Office" palette. (RGB values below are expressed in the
#    BGR-packed long form the PowerPoint object model expects for
#    ColorFormat/ThemeColor.RGB.)
$officeColors = @(
    0,          # dk1     000000
    16777215,   # lt1     FFFFFF
    6968388,    # dk2     44546A
    15132391,   # lt2     E7E6E6
    13998939,   # accent1 5B9BD5
    3243501,    # accent2 ED7D31
    10855845,   # accent3 A5A5A5
    49407,      # accent4 FFC000
    12874308,   # accent5 4472C4
    4697456,    # accent6 70AD47
    12673797,   # hlink   0563C1
    7491477     # folHlink 954F72
)

$themeColors = $s5.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
